$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text number format for cells whose new values
# look like numbers, so they stay text (matching the source data).
# (Looping per-cell since comma multi-area Range strings are not
# reliably applied across all areas by this runtime.)
foreach ($addr in @('D4', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D18', 'D19', 'D20', 'D21', 'D23', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '23.902.48'
$ws.Range('E2').Value = '  -2.82%  '

$ws.Range('D3').Value = '1.623.31'
$ws.Range('E3').Value = '  -2.94%  '

$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.27%  '

$ws.Range('E5').Value = '  +0.36%  '

$ws.Range('D6').Value = '307.48'
$ws.Range('E6').Value = '  -2.34%  '

$ws.Range('D7').Value = '0.3903'
$ws.Range('E7').Value = '  -0.23%  '

$ws.Range('D8').Value = '0.3810'
$ws.Range('E8').Value = '  -3.11%  '

$ws.Range('D9').Value = '1.007'
$ws.Range('E9').Value = '  +0.26%  '

$ws.Range('D10').Value = '49.86'
$ws.Range('E10').Value = '  -4.11%  '

$ws.Range('D11').Value = '1.352'
$ws.Range('E11').Value = '  -2.38%  '

$ws.Range('D12').Value = '0.08446'
$ws.Range('E12').Value = '  -2.38%  '

$ws.Range('D13').Value = '23.71'
$ws.Range('E13').Value = '  -5.46%  '

$ws.Range('D14').Value = '6.972'
$ws.Range('E14').Value = '  -4.38%  '

$ws.Range('D15').Value = '0.00001270'
$ws.Range('E15').Value = '  -3.04%  '

$ws.Range('D16').Value = '7.408'
$ws.Range('E16').Value = '  -4.30%  '

$ws.Range('D17').Value = '1.630.80'
$ws.Range('E17').Value = '  -4.24%  '

$ws.Range('D18').Value = '92.81'
$ws.Range('E18').Value = '  -0.76%  '

$ws.Range('D19').Value = '0.06893'
$ws.Range('E19').Value = '  -2.25%  '

$ws.Range('D20').Value = '19.91'
$ws.Range('E20').Value = '  -3.07%  '

$ws.Range('D21').Value = '6.832'
$ws.Range('E21').Value = '  -3.17%  '

$ws.Range('E22').Value = '  +0.06%  '

$ws.Range('D23').Value = '13.37'
$ws.Range('E23').Value = '  -4.00%  '

$ws.Range('D24').Value = '23.921.78'
$ws.Range('E24').Value = '  -2.79%  '

$ws.Range('D25').Value = '2.387'
$ws.Range('E25').Value = '  +1.50%  '

$ws.Range('D26').Value = '2.856'
$ws.Range('E26').Value = '  +5.18%  '

$ws.Range('D27').Value = '21.99'
$ws.Range('E27').Value = '  -5.21%  '

$ws.Range('D28').Value = '157.34'
$ws.Range('E28').Value = '  -3.22%  '

$ws.Range('D29').Value = '138.19'
$ws.Range('E29').Value = '  -6.13%  '

$ws.Range('D30').Value = '5.245'
$ws.Range('E30').Value = '  -8.53%  '

$ws.Range('D31').Value = '7.717'
$ws.Range('E31').Value = '  -1.59%  '

$ws.Range('E32').Value = '  -2.10%  '

$ws.Range('D33').Value = '1.805.58'
$ws.Range('E33').Value = '  -2.99%  '

$ws.Range('D34').Value = '0.07939'
$ws.Range('E34').Value = '  -5.16%  '

$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.02880'
$ws.Range('E35').Value = '  -4.97%  '

$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').Value = '6.620'
$ws.Range('E36').Value = '  -3.50%  '

$ws.Range('D37').Value = '0.9488'
$ws.Range('E37').Value = '  -2.81%  '

$ws.Range('D38').Value = '0.2644'
$ws.Range('E38').Value = '  -5.64%  '

$ws.Range('D39').Value = '0.09148'
$ws.Range('E39').Value = '  -3.34%  '

$ws.Range('D40').Value = '10.22'
$ws.Range('E40').Value = '  -2.80%  '

$ws.Range('D41').Value = '1.413'
$ws.Range('E41').Value = '  -8.55%  '

$ws.Range('D42').Value = '0.7427'
$ws.Range('E42').Value = '  -5.89%  '

$ws.Range('D43').Value = '13.08'
$ws.Range('E43').Value = '  -2.82%  '

$ws.Range('D44').Value = '15.90'
$ws.Range('E44').Value = '  -2.95%  '

$ws.Range('D45').Value = '0.6802'

$ws.Range('D46').Value = '2.431'
$ws.Range('E46').Value = '  -4.91%  '

$ws.Range('D47').Value = '4.081'
$ws.Range('E47').Value = '  -2.82%  '

$ws.Range('D48').Value = '1.006'
$ws.Range('E48').Value = '  +0.26%  '

$ws.Range('D49').Value = '0.08224'
$ws.Range('E49').Value = '  -4.87%  '

$ws.Range('D50').Value = '132.84'
$ws.Range('E50').Value = '  -3.04%  '

$ws.Range('D51').Value = '1.246'
$ws.Range('E51').Value = '  -5.51%  '
